$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns B through H entirely (their headers/data/widths all go away)
$ws.Range("B1:H2").EntireColumn.Delete()

# Update A1 to the value that used to be in F1 ("input_Name"), keeping its existing (header) style
$ws.Range("A1").Value = "input_Name"

# Clear the value that used to be in A2 ("Codeless Automation Tool") but keep the cell present
$ws.Range("A2").Value = ""
$ws.Range("A2").Font.Bold = $false

# Column A now has the width that used to belong to column F (stored width 12)
$ws.Columns.Item(1).ColumnWidth = 11.17
